# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# The "Estado de Cuenta" worksheet listed 5 worker/period rows (16-20).
# The new data set only has a single worker/period, so:
#   - rows 17-20 (the old extra worker/period entries) are removed
#   - the remaining row 16 is updated to the new worker name / period
#   - the summary counts (Cant. Trabajadores / Cant. Periodos) and the
#     VALOR MORA total are updated to match the single remaining record

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the four obsolete worker/period rows -------------------------
$ws.Rows("17:20").Delete()

# --- update the single remaining worker/period row (row 16) --------------
$ws.Range("D16").Value = "WILMER ARENAS JULIO"
$ws.Range("E16").Value = "2509"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 908526

# --- update the header summary figures ------------------------------------
$ws.Range("E11").Value = 36341
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# --- column D ("Nombre Trabajador") shrinks now that the longest name ----
# --- ("SHEYLA PAOLA AVILA PUELLO") is gone; match the new best-fit width --
$ws.Columns("D:D").ColumnWidth = 20.8
